$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 5
$ws.Cells.Item(5, 8).Value = 110.35714
$ws.Cells.Item(5, 9).Value = 103.75
$ws.Cells.Item(5, 10).Value = 150
$ws.Cells.Item(5, 11).Value = 103.75
$ws.Cells.Item(5, 12).Value = 150
$ws.Cells.Item(5, 13).Value = 11.25
$ws.Cells.Item(5, 14).Value = -380
# Row 40
$ws.Cells.Item(40, 8).Value = 52633056
$ws.Cells.Item(40, 9).Value = 1526.9231
$ws.Cells.Item(40, 10).Value = 166668030
$ws.Cells.Item(40, 11).Value = 1526.9231
$ws.Cells.Item(40, 12).Value = 166668030
$ws.Cells.Item(40, 13).Value = -1351.9231
$ws.Cells.Item(40, 14).Value = -166668380
# Row 98
$ws.Cells.Item(98, 8).Value = 1301.9584
$ws.Cells.Item(98, 9).Value = 1357.5
$ws.Cells.Item(98, 10).Value = 1135.3334
$ws.Cells.Item(98, 11).Value = 1357.5
$ws.Cells.Item(98, 12).Value = 1135.3334
$ws.Cells.Item(98, 13).Value = 140.5
$ws.Cells.Item(98, 14).Value = -4131.3334
# Row 111
$ws.Cells.Item(111, 8).Value = 3465.6667
$ws.Cells.Item(111, 9).Value = 2839.7
$ws.Cells.Item(111, 10).Value = 4717.6
$ws.Cells.Item(111, 11).Value = 8519.099999999999
$ws.Cells.Item(111, 12).Value = 14152.8
$ws.Cells.Item(111, 13).Value = -5452.099999999999
$ws.Cells.Item(111, 14).Value = -20286.8
# Row 115
$ws.Cells.Item(115, 8).Value = 2293.4119
$ws.Cells.Item(115, 9).Value = 698.5
$ws.Cells.Item(115, 10).Value = 3711.111
$ws.Cells.Item(115, 11).Value = 2095.5
$ws.Cells.Item(115, 12).Value = 11133.333
$ws.Cells.Item(115, 13).Value = -528.5
$ws.Cells.Item(115, 14).Value = -14267.333
# Row 122
$ws.Cells.Item(122, 8).Value = 1301.9584
$ws.Cells.Item(122, 9).Value = 1357.5
$ws.Cells.Item(122, 10).Value = 1135.3334
$ws.Cells.Item(122, 11).Value = 4072.5
$ws.Cells.Item(122, 12).Value = 3406.0002
$ws.Cells.Item(122, 13).Value = -1622.5
$ws.Cells.Item(122, 14).Value = -8306.0002
# Row 135
$ws.Cells.Item(135, 8).Value = 296299.72
$ws.Cells.Item(135, 9).Value = 335283.66
$ws.Cells.Item(135, 10).Value = 3920
$ws.Cells.Item(135, 11).Value = 3017552.94
$ws.Cells.Item(135, 12).Value = 35280
$ws.Cells.Item(135, 13).Value = -3015017.94
$ws.Cells.Item(135, 14).Value = -40350
# Row 137
$ws.Cells.Item(137, 8).Value = 2225.6765
$ws.Cells.Item(137, 9).Value = 1257.0588
$ws.Cells.Item(137, 10).Value = 3194.2942
$ws.Cells.Item(137, 11).Value = 3771.1764
$ws.Cells.Item(137, 12).Value = 9582.882599999999
$ws.Cells.Item(137, 13).Value = -1221.1764
$ws.Cells.Item(137, 14).Value = -14682.8826
# Row 138
$ws.Cells.Item(138, 8).Value = 2956.63
$ws.Cells.Item(138, 9).Value = 1558.4138
$ws.Cells.Item(138, 10).Value = 3527.7324
$ws.Cells.Item(138, 11).Value = 4675.2414
$ws.Cells.Item(138, 12).Value = 10583.1972
$ws.Cells.Item(138, 13).Value = 464.7586000000001
$ws.Cells.Item(138, 14).Value = -20863.1972
# Row 141
$ws.Cells.Item(141, 8).Value = 2308.8481
$ws.Cells.Item(141, 9).Value = 1911.5
$ws.Cells.Item(141, 10).Value = 2503.7737
$ws.Cells.Item(141, 11).Value = 5734.5
$ws.Cells.Item(141, 12).Value = 7511.321100000001
$ws.Cells.Item(141, 13).Value = -554.5
$ws.Cells.Item(141, 14).Value = -17871.3211

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Cells.Item(32, 8).Value = 27706.545
$ws.Cells.Item(32, 9).Value = 10427.712
$ws.Cells.Item(32, 10).Value = 131379.55
$ws.Cells.Item(32, 11).Value = 10427.712
$ws.Cells.Item(32, 12).Value = 131379.55
$ws.Cells.Item(32, 13).Value = -10140.712
$ws.Cells.Item(32, 14).Value = -131953.55
# Row 61
$ws.Cells.Item(61, 8).Value = 2065.8
$ws.Cells.Item(61, 9).Value = 1984
$ws.Cells.Item(61, 10).Value = 2597.5
$ws.Cells.Item(61, 11).Value = 1984
$ws.Cells.Item(61, 12).Value = 2597.5
$ws.Cells.Item(61, 13).Value = -1772
$ws.Cells.Item(61, 14).Value = -3021.5
# Row 63
$ws.Cells.Item(63, 8).Value = 2425.4348
$ws.Cells.Item(63, 9).Value = 1861.5625
$ws.Cells.Item(63, 10).Value = 3714.2856
$ws.Cells.Item(63, 11).Value = 1861.5625
$ws.Cells.Item(63, 12).Value = 3714.2856
$ws.Cells.Item(63, 13).Value = -1175.5625
$ws.Cells.Item(63, 14).Value = -5086.2856
# Row 66
$ws.Cells.Item(66, 8).Value = 2425.4348
$ws.Cells.Item(66, 9).Value = 1861.5625
$ws.Cells.Item(66, 10).Value = 3714.2856
$ws.Cells.Item(66, 11).Value = 9307.8125
$ws.Cells.Item(66, 12).Value = 18571.428
$ws.Cells.Item(66, 13).Value = -5875.8125
$ws.Cells.Item(66, 14).Value = -25435.428
# Row 74
$ws.Cells.Item(74, 8).Value = 2929.6223
$ws.Cells.Item(74, 9).Value = 2661.8215
$ws.Cells.Item(74, 10).Value = 3370.7058
$ws.Cells.Item(74, 11).Value = 2661.8215
$ws.Cells.Item(74, 12).Value = 3370.7058
$ws.Cells.Item(74, 13).Value = -1787.8215
$ws.Cells.Item(74, 14).Value = -5118.7058
# Row 77
$ws.Cells.Item(77, 8).Value = 2929.6223
$ws.Cells.Item(77, 9).Value = 2661.8215
$ws.Cells.Item(77, 10).Value = 3370.7058
$ws.Cells.Item(77, 11).Value = 13309.1075
$ws.Cells.Item(77, 12).Value = 16853.529
$ws.Cells.Item(77, 13).Value = -8941.1075
$ws.Cells.Item(77, 14).Value = -25589.529
# Row 97
$ws.Cells.Item(97, 8).Value = 1030.2941
$ws.Cells.Item(97, 9).Value = 651.25
$ws.Cells.Item(97, 10).Value = 1940
$ws.Cells.Item(97, 11).Value = 651.25
$ws.Cells.Item(97, 12).Value = 1940
$ws.Cells.Item(97, 13).Value = -155.25
$ws.Cells.Item(97, 14).Value = -2932
# Row 132
$ws.Cells.Item(132, 8).Value = 3058.25
$ws.Cells.Item(132, 9).Value = 1860.1786
$ws.Cells.Item(132, 10).Value = 7251.5
$ws.Cells.Item(132, 11).Value = 5580.5358
$ws.Cells.Item(132, 12).Value = 21754.5
$ws.Cells.Item(132, 13).Value = -3050.5358
$ws.Cells.Item(132, 14).Value = -26814.5
# Row 133
$ws.Cells.Item(133, 8).Value = 46500
$ws.Cells.Item(133, 9).Value = 0
$ws.Cells.Item(133, 10).Value = 46500
$ws.Cells.Item(133, 11).Value = 0
$ws.Cells.Item(133, 12).Value = 46500
$ws.Cells.Item(133, 14).Value = -51560
# Row 136
$ws.Cells.Item(136, 8).Value = 2065.8
$ws.Cells.Item(136, 9).Value = 1984
$ws.Cells.Item(136, 10).Value = 2597.5
$ws.Cells.Item(136, 11).Value = 5952
$ws.Cells.Item(136, 12).Value = 7792.5
$ws.Cells.Item(136, 13).Value = -3402
$ws.Cells.Item(136, 14).Value = -12892.5

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 94
$ws.Cells.Item(94, 8).Value = 575.1667
$ws.Cells.Item(94, 9).Value = 493.96295
$ws.Cells.Item(94, 10).Value = 721.3333
$ws.Cells.Item(94, 11).Value = 493.96295
$ws.Cells.Item(94, 12).Value = 721.3333
$ws.Cells.Item(94, 13).Value = -42.96294999999998
$ws.Cells.Item(94, 14).Value = -1623.3333
# Row 134
$ws.Cells.Item(134, 8).Value = 1892.3438
$ws.Cells.Item(134, 9).Value = 1488.3043
$ws.Cells.Item(134, 10).Value = 2924.889
$ws.Cells.Item(134, 11).Value = 4464.9129
$ws.Cells.Item(134, 12).Value = 8774.667000000001
$ws.Cells.Item(134, 13).Value = -1929.9129
$ws.Cells.Item(134, 14).Value = -13844.667
# Row 141
$ws.Cells.Item(141, 8).Value = 78161.71000000001
$ws.Cells.Item(141, 9).Value = 0
$ws.Cells.Item(141, 10).Value = 78161.71000000001
$ws.Cells.Item(141, 11).Value = 0
$ws.Cells.Item(141, 12).Value = 78161.71000000001
$ws.Cells.Item(141, 14).Value = -88521.71000000001

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Cells.Item(31, 8).Value = 3675.8823
$ws.Cells.Item(31, 9).Value = 1815.1515
$ws.Cells.Item(31, 10).Value = 7087.222
$ws.Cells.Item(31, 11).Value = 1815.1515
$ws.Cells.Item(31, 12).Value = 7087.222
$ws.Cells.Item(31, 13).Value = -1520.1515
$ws.Cells.Item(31, 14).Value = -7677.222
# Row 34
$ws.Cells.Item(34, 8).Value = 3675.8823
$ws.Cells.Item(34, 9).Value = 1815.1515
$ws.Cells.Item(34, 10).Value = 7087.222
$ws.Cells.Item(34, 11).Value = 1815.1515
$ws.Cells.Item(34, 12).Value = 7087.222
$ws.Cells.Item(34, 13).Value = -1613.1515
$ws.Cells.Item(34, 14).Value = -7491.222
# Row 58
$ws.Cells.Item(58, 8).Value = 1531.2609
$ws.Cells.Item(58, 9).Value = 1635.95
$ws.Cells.Item(58, 10).Value = 833.3333
$ws.Cells.Item(58, 11).Value = 1635.95
$ws.Cells.Item(58, 12).Value = 833.3333
$ws.Cells.Item(58, 13).Value = -1432.95
$ws.Cells.Item(58, 14).Value = -1239.3333
# Row 99
$ws.Cells.Item(99, 8).Value = 1544.9412
$ws.Cells.Item(99, 9).Value = 1412.5
$ws.Cells.Item(99, 10).Value = 1662.6666
$ws.Cells.Item(99, 11).Value = 1412.5
$ws.Cells.Item(99, 12).Value = 1662.6666
$ws.Cells.Item(99, 13).Value = 85.5
$ws.Cells.Item(99, 14).Value = -4658.6666
# Row 126
$ws.Cells.Item(126, 8).Value = 1544.9412
$ws.Cells.Item(126, 9).Value = 1412.5
$ws.Cells.Item(126, 10).Value = 1662.6666
$ws.Cells.Item(126, 11).Value = 4237.5
$ws.Cells.Item(126, 12).Value = 4987.9998
$ws.Cells.Item(126, 13).Value = -1767.5
$ws.Cells.Item(126, 14).Value = -9927.9998
# Row 132
$ws.Cells.Item(132, 8).Value = 2368.9678
$ws.Cells.Item(132, 9).Value = 1889.125
$ws.Cells.Item(132, 10).Value = 4014.1428
$ws.Cells.Item(132, 11).Value = 5667.375
$ws.Cells.Item(132, 12).Value = 12042.4284
$ws.Cells.Item(132, 13).Value = -3137.375
$ws.Cells.Item(132, 14).Value = -17102.4284
# Row 134
$ws.Cells.Item(134, 8).Value = 5442.3335
$ws.Cells.Item(134, 9).Value = 6167.087
$ws.Cells.Item(134, 10).Value = 1275
$ws.Cells.Item(134, 11).Value = 18501.261
$ws.Cells.Item(134, 12).Value = 3825
$ws.Cells.Item(134, 13).Value = -15966.261
$ws.Cells.Item(134, 14).Value = -8895
# Row 136
$ws.Cells.Item(136, 8).Value = 1531.2609
$ws.Cells.Item(136, 9).Value = 1635.95
$ws.Cells.Item(136, 10).Value = 833.3333
$ws.Cells.Item(136, 11).Value = 4907.85
$ws.Cells.Item(136, 12).Value = 2499.9999
$ws.Cells.Item(136, 13).Value = -2357.85
$ws.Cells.Item(136, 14).Value = -7599.9999

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 113
$ws.Cells.Item(113, 8).Value = 689.125
$ws.Cells.Item(113, 9).Value = 661.75
$ws.Cells.Item(113, 10).Value = 826
$ws.Cells.Item(113, 11).Value = 1985.25
$ws.Cells.Item(113, 12).Value = 2478
$ws.Cells.Item(113, 13).Value = 184.75
$ws.Cells.Item(113, 14).Value = -6818
# Row 132
$ws.Cells.Item(132, 8).Value = 576643.75
$ws.Cells.Item(132, 9).Value = 941236.1
$ws.Cells.Item(132, 10).Value = 9500
$ws.Cells.Item(132, 11).Value = 8471124.9
$ws.Cells.Item(132, 12).Value = 85500
$ws.Cells.Item(132, 13).Value = -8468594.9
$ws.Cells.Item(132, 14).Value = -90560

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 122
$ws.Cells.Item(122, 8).Value = 2210.2942
$ws.Cells.Item(122, 9).Value = 2148.2144
$ws.Cells.Item(122, 10).Value = 2500
$ws.Cells.Item(122, 11).Value = 6444.6432
$ws.Cells.Item(122, 12).Value = 7500
$ws.Cells.Item(122, 13).Value = -3994.6432
$ws.Cells.Item(122, 14).Value = -12400
# Row 132
$ws.Cells.Item(132, 8).Value = 2832.1333
$ws.Cells.Item(132, 9).Value = 2806.3076
$ws.Cells.Item(132, 10).Value = 3000
$ws.Cells.Item(132, 11).Value = 8418.9228
$ws.Cells.Item(132, 12).Value = 9000
$ws.Cells.Item(132, 13).Value = -5888.9228
$ws.Cells.Item(132, 14).Value = -14060

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 40
$ws.Cells.Item(40, 8).Value = 4201
$ws.Cells.Item(40, 9).Value = 0
$ws.Cells.Item(40, 10).Value = 4201
$ws.Cells.Item(40, 11).Value = 0
$ws.Cells.Item(40, 12).Value = 4201
$ws.Cells.Item(40, 14).Value = -4473
$ws.Cells.Item(40, 13).ClearContents()
# Row 46
$ws.Cells.Item(46, 9).Value = 1342.8572
$ws.Cells.Item(46, 10).Value = 4333.3335
$ws.Cells.Item(46, 11).Value = 1342.8572
$ws.Cells.Item(46, 12).Value = 4333.3335
$ws.Cells.Item(46, 13).Value = -1154.8572
$ws.Cells.Item(46, 14).Value = -4709.3335
# Row 61
$ws.Cells.Item(61, 8).Value = 3681.3635
$ws.Cells.Item(61, 9).Value = 2249.1667
$ws.Cells.Item(61, 10).Value = 5400
$ws.Cells.Item(61, 11).Value = 2249.1667
$ws.Cells.Item(61, 12).Value = 5400
$ws.Cells.Item(61, 13).Value = -2047.1667
$ws.Cells.Item(61, 14).Value = -5804
# Row 93
$ws.Cells.Item(93, 8).Value = 1356.7046
$ws.Cells.Item(93, 9).Value = 1178.9354
$ws.Cells.Item(93, 10).Value = 1780.6154
$ws.Cells.Item(93, 11).Value = 1178.9354
$ws.Cells.Item(93, 12).Value = 1780.6154
$ws.Cells.Item(93, 13).Value = 69.06459999999993
$ws.Cells.Item(93, 14).Value = -4276.6154
# Row 113
$ws.Cells.Item(113, 8).Value = 3681.3635
$ws.Cells.Item(113, 9).Value = 2249.1667
$ws.Cells.Item(113, 10).Value = 5400
$ws.Cells.Item(113, 11).Value = 2249.1667
$ws.Cells.Item(113, 12).Value = 5400
$ws.Cells.Item(113, 13).Value = -79.16670000000022
$ws.Cells.Item(113, 14).Value = -9740
# Row 122
$ws.Cells.Item(122, 8).Value = 4799
$ws.Cells.Item(122, 9).Value = 4651
$ws.Cells.Item(122, 10).Value = 4848.3335
$ws.Cells.Item(122, 11).Value = 13953
$ws.Cells.Item(122, 12).Value = 14545.0005
$ws.Cells.Item(122, 13).Value = -11503
$ws.Cells.Item(122, 14).Value = -19445.0005
# Row 136
$ws.Cells.Item(136, 8).Value = 2006.3529
$ws.Cells.Item(136, 9).Value = 1819.25
$ws.Cells.Item(136, 10).Value = 5000
$ws.Cells.Item(136, 11).Value = 5457.75
$ws.Cells.Item(136, 12).Value = 15000
$ws.Cells.Item(136, 13).Value = -2907.75
$ws.Cells.Item(136, 14).Value = -20100

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 100
$ws.Cells.Item(100, 8).Value = 909891.4399999999
$ws.Cells.Item(100, 9).Value = 485.7143
$ws.Cells.Item(100, 10).Value = 2501351.5
$ws.Cells.Item(100, 11).Value = 971.4286
$ws.Cells.Item(100, 12).Value = 5002703
$ws.Cells.Item(100, 13).Value = -430.4286
$ws.Cells.Item(100, 14).Value = -5003785
# Row 132
$ws.Cells.Item(132, 8).Value = 2903.4211
$ws.Cells.Item(132, 9).Value = 4908.643
$ws.Cells.Item(132, 10).Value = 1733.7084
$ws.Cells.Item(132, 11).Value = 14725.929
$ws.Cells.Item(132, 12).Value = 5201.1252
$ws.Cells.Item(132, 13).Value = -12195.929
$ws.Cells.Item(132, 14).Value = -10261.1252
